$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.084.61'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '2.585.19'
$ws.Range("E3").Value = '  +8.76%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.89'
$ws.Range("E5").Value = '  +1.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.82'
$ws.Range("E6").Value = '  +1.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.592'
$ws.Range("E7").Value = '  +4.61%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.572'
$ws.Range("E9").Value = '  +12.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.41'
$ws.Range("E10").Value = '  +11.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0836'
$ws.Range("E11").Value = '  +6.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.12'
$ws.Range("E12").Value = '  +14.27%  '
$ws.Range("D13").Value = '2.973.92'
$ws.Range("E13").Value = '  +8.39%  '
$ws.Range("E14").Value = '  +1.20%  '
$ws.Range("D15").Value = '2.581.09'
$ws.Range("E15").Value = '  +8.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.896'
$ws.Range("E16").Value = '  +8.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.76'
$ws.Range("E17").Value = '  +7.44%  '
$ws.Range("D18").Value = '46.141.87'
$ws.Range("E18").Value = '  +0.52%  '
$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.04'
$ws.Range("E19").Value = '  +3.25%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000101'
$ws.Range("E20").Value = '  +6.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.62'
$ws.Range("E21").Value = '  +9.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.77'
$ws.Range("E22").Value = '  +6.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '252.94'
$ws.Range("E23").Value = '  +4.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.98'
$ws.Range("E24").Value = '  +6.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.19'
$ws.Range("E25").Value = '  +13.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.80'
$ws.Range("E26").Value = '  +33.22%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.38'
$ws.Range("E28").Value = '  +6.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.32'
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.26'
$ws.Range("E30").Value = '  +2.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.06'
$ws.Range("E31").Value = '  +9.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.68'
$ws.Range("E32").Value = '  -1.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.93'
$ws.Range("E33").Value = '  +4.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.26'
$ws.Range("E34").Value = '  +18.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '152.05'
$ws.Range("E35").Value = '  +3.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0821'
$ws.Range("E36").Value = '  +6.59%  '
$ws.Range("E37").Value = '  +2.39%  '
$ws.Range("E38").Value = '  +4.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.18'
$ws.Range("E39").Value = '  +8.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.16'
$ws.Range("E40").Value = '  +7.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.55'
$ws.Range("E41").Value = '  +11.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0319'
$ws.Range("E42").Value = '  +6.99%  '
$ws.Range("D43").Value = '2.060.70'
$ws.Range("E43").Value = '  +6.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.61'
$ws.Range("E44").Value = '  +38.55%  '
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.10'
$ws.Range("E46").Value = '  -1.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.22'
$ws.Range("E47").Value = '  +8.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.77'
$ws.Range("E48").Value = '  -1.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '108.32'
$ws.Range("E49").Value = '  +9.52%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.200'
$ws.Range("E50").Value = '  +7.81%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.834.79'
$ws.Range("E51").Value = '  +8.38%  '
